$d = $word.ActiveDocument

function Replace-Label {
    param(
        [string]$OldText,
        [string]$NewText
    )

    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Text = $OldText
    $find.Replacement.Text = $NewText
    $find.Replacement.Font.Bold = $true
    $find.Forward = $true
    $find.Wrap = 1
    $find.Format = $true
    $find.MatchCase = $true
    $find.MatchWholeWord = $false
    $find.MatchWildcards = $false
    $find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $true, $NewText, 2)
}

Replace-Label 'Chiffre d’affaires (K€)' 'Revenu (K$)'
Replace-Label 'Coût des marchandises vendues ($K)' 'Coût des marchandises vendues (K$)'
Replace-Label 'Marge bénéficiaire brute (%)' 'Marge bénéficiaire brut (%)'
Replace-Label 'Dépenses de fonctionnement ($K)' 'Dépenses opérationnelles (K$)'
Replace-Label 'EBITDA ($K)' 'EBITDA (K$)'
Replace-Label 'Charges d’intérêt ($K)' 'Charges d’intérêts (K$)'
Replace-Label 'Bénéfice avant impôts ($K)' 'Bénéfice avant impôt (K$)'
Replace-Label 'Revenus nets ($K)' 'Résultat net (K$)'
Replace-Label 'Total des actifs ($K)' 'Total actif (K$)'
Replace-Label 'Total du passif ($K)' 'Total passif (K$)'
Replace-Label 'Capitaux propres ($K)' 'Capitaux propres des actionnaires (K$)'
